# Generate Report for handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) values in row 5 of the
# zh-cn and de-de sheets to reflect a newly generated report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-26 10:08:48"
$wsZhCn.Range("G5").Value = "2016-01-26 10:09:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-26 10:08:59"
$wsDeDe.Range("G5").Value = "2016-01-26 10:09:52"
